$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix comma-separated names to use periods (per commit: fixed formatting) ---
$ws.Range("E89").Value = "RICCOTTI. MARIANA EDITH"
$ws.Range("F96").Value = "MERCANZINI. GASTON ARIEL"
$ws.Range("E114").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E126").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

# --- Fix floating point 'Importe' values: was scraped with Spanish formatting
# (period=thousands, comma=decimal); now stored as plain 1234.56 with a dot decimal ---
$importeFixes = @(
    @{Addr="H2"; New="8990.50"},
    @{Addr="H3"; New="14028.31"},
    @{Addr="H4"; New="4840.00"},
    @{Addr="H5"; New="432.82"},
    @{Addr="H6"; New="200.00"},
    @{Addr="H7"; New="3600.00"},
    @{Addr="H8"; New="442824.62"},
    @{Addr="H9"; New="5426.00"},
    @{Addr="H10"; New="52150.04"},
    @{Addr="H11"; New="15747.00"},
    @{Addr="H12"; New="20226.88"},
    @{Addr="H13"; New="2203.00"},
    @{Addr="H14"; New="18675.16"},
    @{Addr="H15"; New="19639.46"},
    @{Addr="H16"; New="382.50"},
    @{Addr="H17"; New="7340.00"},
    @{Addr="H18"; New="745.00"},
    @{Addr="H19"; New="5300.00"},
    @{Addr="H20"; New="12087.67"},
    @{Addr="H21"; New="1706.05"},
    @{Addr="H22"; New="4529.37"},
    @{Addr="H23"; New="441.60"},
    @{Addr="H24"; New="1155.50"},
    @{Addr="H25"; New="3715.53"},
    @{Addr="H26"; New="21944.80"},
    @{Addr="H27"; New="8775.00"},
    @{Addr="H28"; New="76020.90"},
    @{Addr="H29"; New="139292.25"},
    @{Addr="H30"; New="4598.00"},
    @{Addr="H31"; New="25203.99"},
    @{Addr="H32"; New="2761.42"},
    @{Addr="H33"; New="4236.00"},
    @{Addr="H34"; New="11062.16"},
    @{Addr="H35"; New="1229.28"},
    @{Addr="H36"; New="2080.00"},
    @{Addr="H37"; New="864.72"},
    @{Addr="H38"; New="246.27"},
    @{Addr="H39"; New="36998.90"},
    @{Addr="H40"; New="230.00"},
    @{Addr="H41"; New="190227.31"},
    @{Addr="H42"; New="15870.00"},
    @{Addr="H43"; New="5500.00"},
    @{Addr="H44"; New="951.76"},
    @{Addr="H45"; New="413.20"},
    @{Addr="H46"; New="1481.72"},
    @{Addr="H47"; New="2348.00"},
    @{Addr="H48"; New="19000.00"},
    @{Addr="H49"; New="1638.30"},
    @{Addr="H50"; New="90.00"},
    @{Addr="H51"; New="17283.00"},
    @{Addr="H52"; New="1610.00"},
    @{Addr="H53"; New="2400.00"},
    @{Addr="H54"; New="81127.04"},
    @{Addr="H55"; New="1991.25"},
    @{Addr="H56"; New="600.00"},
    @{Addr="H57"; New="165.00"},
    @{Addr="H58"; New="3600.00"},
    @{Addr="H59"; New="400.00"},
    @{Addr="H60"; New="190705.00"},
    @{Addr="H61"; New="20030.00"},
    @{Addr="H62"; New="6665.00"},
    @{Addr="H63"; New="1720.00"},
    @{Addr="H64"; New="8480.00"},
    @{Addr="H65"; New="238.00"},
    @{Addr="H66"; New="372.60"},
    @{Addr="H67"; New="2906.00"},
    @{Addr="H68"; New="5600.00"},
    @{Addr="H69"; New="391178.31"},
    @{Addr="H70"; New="53261.49"},
    @{Addr="H71"; New="11.27"},
    @{Addr="H72"; New="7800.00"},
    @{Addr="H73"; New="629.32"},
    @{Addr="H74"; New="120.00"},
    @{Addr="H75"; New="21624.00"},
    @{Addr="H76"; New="2434.00"},
    @{Addr="H77"; New="4122.50"},
    @{Addr="H78"; New="970.00"},
    @{Addr="H79"; New="1760.00"},
    @{Addr="H80"; New="27888.00"},
    @{Addr="H81"; New="114.00"},
    @{Addr="H82"; New="1815.00"},
    @{Addr="H83"; New="340.00"},
    @{Addr="H84"; New="33297.90"},
    @{Addr="H85"; New="3439.50"},
    @{Addr="H86"; New="778.00"},
    @{Addr="H87"; New="423.86"},
    @{Addr="H88"; New="568.70"},
    @{Addr="H89"; New="2000.00"},
    @{Addr="H90"; New="14800.00"},
    @{Addr="H91"; New="7949.00"},
    @{Addr="H92"; New="1806.00"},
    @{Addr="H93"; New="2500.00"},
    @{Addr="H94"; New="1300.00"},
    @{Addr="H95"; New="4000.00"},
    @{Addr="H96"; New="9000.00"},
    @{Addr="H97"; New="3080.00"},
    @{Addr="H98"; New="456.50"},
    @{Addr="H99"; New="3505.06"},
    @{Addr="H100"; New="1000.00"},
    @{Addr="H101"; New="4000.00"},
    @{Addr="H102"; New="1657.50"},
    @{Addr="H103"; New="19199.98"},
    @{Addr="H104"; New="800.00"},
    @{Addr="H105"; New="1548.86"},
    @{Addr="H106"; New="3600.00"},
    @{Addr="H107"; New="1196.00"},
    @{Addr="H108"; New="1030.00"},
    @{Addr="H109"; New="22500.00"},
    @{Addr="H110"; New="726.00"},
    @{Addr="H111"; New="1700.00"},
    @{Addr="H112"; New="780.00"},
    @{Addr="H113"; New="279.11"},
    @{Addr="H114"; New="7700.00"},
    @{Addr="H115"; New="35.91"},
    @{Addr="H116"; New="1410.00"},
    @{Addr="H117"; New="2129.60"},
    @{Addr="H118"; New="164.74"},
    @{Addr="H119"; New="12200.00"},
    @{Addr="H120"; New="2005.24"},
    @{Addr="H121"; New="844.00"},
    @{Addr="H122"; New="2796.00"},
    @{Addr="H123"; New="9910.00"},
    @{Addr="H124"; New="2840.79"},
    @{Addr="H125"; New="14485.00"},
    @{Addr="H126"; New="19780.00"},
    @{Addr="H127"; New="314.38"},
    @{Addr="H128"; New="4816.82"},
    @{Addr="H129"; New="1122.22"},
    @{Addr="H130"; New="90.00"},
    @{Addr="H131"; New="3600.00"},
    @{Addr="H132"; New="95182.50"},
    @{Addr="H133"; New="430.00"},
    @{Addr="H134"; New="3860.00"},
    @{Addr="H135"; New="2120.00"},
    @{Addr="H136"; New="1400.00"},
    @{Addr="H137"; New="305.50"},
    @{Addr="H138"; New="8735.23"},
    @{Addr="H139"; New="2607.13"},
    @{Addr="H140"; New="837922.72"},
    @{Addr="H141"; New="2300.00"},
    @{Addr="H142"; New="5464.00"},
    @{Addr="H143"; New="101340.00"},
    @{Addr="H144"; New="9870.00"},
    @{Addr="H145"; New="4800.00"},
    @{Addr="H146"; New="960.74"},
    @{Addr="H147"; New="73278.00"},
    @{Addr="H148"; New="6300.00"},
    @{Addr="H149"; New="4500.00"},
    @{Addr="H150"; New="2900.00"},
    @{Addr="H151"; New="6915.00"},
    @{Addr="H152"; New="7070.00"},
    @{Addr="H153"; New="3500.00"}
)

foreach ($fix in $importeFixes) {
    $cell = $ws.Range($fix.Addr)
    $savedStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $fix.New
    $cell.Style = $savedStyle
}
